# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values on the active worksheet to reflect
# the recalculated strikeout (K) counts instead of the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$gUpdates = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 3
    16 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    23 = 1
    25 = 1
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}
